$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Prix Spot": add a new last column AK ("20-jul") mirroring the
# existing day columns (B .. AJ), with header cell AK1 taking on the
# same bold/border/centered style as the other header cells.
# ------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the previous header cell (AJ1, col 36) onto the
# new header cell (AK1, col 37), then overwrite its value/text.
$wsPrix.Cells.Item(1, 36).Copy($wsPrix.Cells.Item(1, 37))
$wsPrix.Cells.Item(1, 37).Value = "20-jul"

$prixValues = @(
    51.61,
    43.01,
    25.5,
    22.64,
    21.83,
    19.99,
    16.32,
    20.76,
    3.6,
    0.65,
    0,
    0,
    0,
    -0.01,
    -0.01,
    0,
    3.08,
    0,
    5.99,
    11.17,
    43.81,
    50.8,
    78.33,
    74.93
)

for ($i = 0; $i -lt $prixValues.Length; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 37).Value = $prixValues[$i]
}

# ------------------------------------------------------------------
# Sheet "Gaz": append a new row 34 for 2025-07-18.
# ------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force column A to be treated as text (not auto-parsed into a date
# serial) while writing, then restore the default "General" format so
# no style index ends up stamped on the new cell (matches existing
# data rows, which carry no explicit style).
$wsGaz.Cells.Item(34, 1).NumberFormat = "@"
$wsGaz.Cells.Item(34, 1).Value = "2025-07-18"
$wsGaz.Cells.Item(34, 1).Style = "Normal"
$wsGaz.Cells.Item(34, 2).Value = 32.85

# ------------------------------------------------------------------
# Sheet "CO2": append a new row 34 for 2025-07-18.
# ------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Cells.Item(34, 1).NumberFormat = "@"
$wsCo2.Cells.Item(34, 1).Value = "2025-07-18"
$wsCo2.Cells.Item(34, 1).Style = "Normal"
$wsCo2.Cells.Item(34, 2).Value = 69.2
